# AMOS B05 - Agile Programming.pptx
# Summer 2024 update: swap the table style used by the table on the
# "practices overview" slide (slide 9) from the old style GUID to the
# new one.

$p = $ppt.ActivePresentation

# The table lives on slide 9 (the graphicFrame named "Google Shape;91;p16").
$s = $p.Slides.Item(9)

# Locate the shape that actually has a table (defensive: don't hardcode
# the shape index in case collection ordering ever shifts).
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{2C25DEFE-EB13-46BB-99A0-56F9EF83F260}")
}
